# Aashish Sort 1 data
# Fill in the previously-empty "Aashish Sort 1" benchmark rows (117-121)
# on Sheet1 with the recorded trial timings for each input size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 117: Trial 1
$ws.Range("P117").Value = 2
$ws.Range("Q117").Value = 4
$ws.Range("R117").Value = 9
$ws.Range("S117").Value = 78
$ws.Range("T117").Value = 1084
$ws.Range("U117").Value = 4684

# Row 118: Trial 2
$ws.Range("P118").Value = 1
$ws.Range("Q118").Value = 3
$ws.Range("R118").Value = 15
$ws.Range("S118").Value = 113
$ws.Range("T118").Value = 961
$ws.Range("U118").Value = 9842

# Row 119: Trial 3
$ws.Range("P119").Value = 1
$ws.Range("Q119").Value = 1
$ws.Range("R119").Value = 8
$ws.Range("S119").Value = 58
$ws.Range("T119").Value = 375
$ws.Range("U119").Value = 9862

# Row 120: Trial 4
$ws.Range("P120").Value = 1
$ws.Range("Q120").Value = 1
$ws.Range("R120").Value = 20
$ws.Range("S120").Value = 46
$ws.Range("T120").Value = 903
$ws.Range("U120").Value = 9694

# Row 121: Trial 5
$ws.Range("P121").Value = 1
$ws.Range("Q121").Value = 1
$ws.Range("R121").Value = 27
$ws.Range("S121").Value = 48
$ws.Range("T121").Value = 916
$ws.Range("U121").Value = 4127

# Reflect the author's final scroll/selection position.
$ws.Range("U122").Select()
